$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9591774940490723
$ws.Range("B1").Value = 2.810628414154053
$ws.Range("C1").Value = 5.46741247177124
$ws.Range("D1").Value = 2.099783420562744
$ws.Range("E1").Value = 1.187044262886047
